# Actualización automática 2025-11-27 08:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": row 6 (RIOS CARRION ANGEL BENIGNO / CERAMICAS AL COSTO S.A.S.) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L6").Value = 739.0700000000001    # PIEDRA SINTERIZADA
$wsGrupo.Range("M6").Value = 21624.44             # PORCELANATO

# --- Sheet "VENTA MENSUAL": row 6 (RIOS CARRION ANGEL BENIGNO / CERAMICAS AL COSTO S.A.S.) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F6").Value = 22363.51           # noviembre

# --- Sheet "VENTA MENSUAL": row 26 (TOTAL) ---
$wsMensual.Range("F26").Value = 56295.38          # noviembre total

# --- Sheet "CUMPLIMIENTO MENSUAL": row 11 (PIEDRA SINTERIZADA) ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D11").Value = 739.0700000000001
$wsCumpl.Range("E11").Value = 163.8099999999999
$wsCumpl.Range("F11").Value = 0.8185694665957824

# --- Sheet "CUMPLIMIENTO MENSUAL": row 12 (PORCELANATO) ---
$wsCumpl.Range("D12").Value = 55556.31
$wsCumpl.Range("E12").Value = -20855.31
$wsCumpl.Range("F12").Value = 1.60100025935852

# --- Sheet "CUMPLIMIENTO MENSUAL": row 14 (TOTAL) ---
$wsCumpl.Range("D14").Value = 56295.38
$wsCumpl.Range("E14").Value = -15517.63941051808
$wsCumpl.Range("F14").Value = 1.380541912970054
